# "Generate Report for handback"
#
# For both the zh-cn and de-de sheets, the two pending files
# (175b4d79-...md and c1775312-...md) have now been handed back:
#   - Status (col B) flips from "Not yet handed off" to "Handed back"
#   - The "Latest Target File" (col E) / "Latest Handback File" (col F)
#     hyperlinks get populated (mirroring the existing Source/Handoff
#     hyperlinks in cols A/C)
#   - "Latest Handback DateTime" (col G) gets a real timestamp
#   - "Handoff Reason" (col H) is "Include" for these rows

$wb = $excel.ActiveWorkbook

function Update-LocSheet($SheetName, $MdUrl1, $XlfUrl1, $MdUrl2, $XlfUrl2, $HandbackTime) {

    $ws = $wb.Worksheets.Item($SheetName)

    $mdText1  = $ws.Range("A2").Value2
    $xlfText1 = $ws.Range("C2").Value2
    $mdText2  = $ws.Range("A3").Value2
    $xlfText2 = $ws.Range("C3").Value2

    # Row 2 (175b4d79-...)
    $ws.Range("B2").Value = "Handed back"
    $ws.Hyperlinks.Add($ws.Range("E2"), $MdUrl1, "", "", $mdText1)
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfUrl1, "", "", $xlfText1)
    $ws.Range("G2").Value = $HandbackTime
    $ws.Range("H2").Value = "Include"

    # Row 3 (c1775312-...)
    $ws.Range("B3").Value = "Handed back"
    $ws.Hyperlinks.Add($ws.Range("E3"), $MdUrl2, "", "", $mdText2)
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfUrl2, "", "", $xlfText2)
    $ws.Range("G3").Value = $HandbackTime
    $ws.Range("H3").Value = "Include"
}

Update-LocSheet "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1e0e097c2febc7957cde2258e87d8e642b2c195d/e2e/175b4d79-3667-47f9-a108-49103b0a3086.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6fcc0e7f9e00ff2679f7e0ead41dbb4ef489797c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/175b4d79-3667-47f9-a108-49103b0a3086.0c32bb012f9794161e918237996a415e371ef29b.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1e0e097c2febc7957cde2258e87d8e642b2c195d/e2e/c1775312-ebdb-49c5-82dc-6312cb165ef6.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6fcc0e7f9e00ff2679f7e0ead41dbb4ef489797c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/c1775312-ebdb-49c5-82dc-6312cb165ef6.aece7cc6cae98eeb6be0fca0cc7739ae6e6bcde6.zh-cn.xlf" `
    "2016-01-08 11:46:25"

Update-LocSheet "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1e0e097c2febc7957cde2258e87d8e642b2c195d/e2e/175b4d79-3667-47f9-a108-49103b0a3086.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d86f29623345b66ea7e02c8e20c1fab2ea92da66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/175b4d79-3667-47f9-a108-49103b0a3086.0c32bb012f9794161e918237996a415e371ef29b.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/1e0e097c2febc7957cde2258e87d8e642b2c195d/e2e/c1775312-ebdb-49c5-82dc-6312cb165ef6.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d86f29623345b66ea7e02c8e20c1fab2ea92da66/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/c1775312-ebdb-49c5-82dc-6312cb165ef6.aece7cc6cae98eeb6be0fca0cc7739ae6e6bcde6.de-de.xlf" `
    "2016-01-08 11:46:46"
